$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (rows 2 through 516): 45178 -> 45179.
$lastRow = 516
$ws.Range("C2:C$lastRow").Value = 45179
